{"js": "// Word JS API (Office.js) script.\n// Body of: async (context) => { ... }\n//\n// Target change (pr1): the paragraph's single run \"yayaysysys\" becomes\n// three runs with identical rPr (<w:lang w:val=\"en-US\"/>):\n//   \"Y\"  +  \"ayaysysys\"  +  \"123123213\"\n// i.e. the leading \"y\" is capitalized to \"Y\" and \"123123213\" is appended\n// at the end, but the text stays split across three separate <w:r>\n// elements (as real Word leaves them after two distinct edit operations)\n// instead of being coalesced back into one run.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// pr1 is the (only) paragraph containing the original text.\nconst target = paragraphs.items[0];\n\n// --- Edit 1: capitalize the leading \"y\" -> \"Y\" -----------------------\nconst contentRange = target.getRange(\"Content\");\nconst hits = contentRange.search(\"y\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nconst firstY = hits.items[0];\nconst runY = firstY.insertText(\"Y\", \"Replace\");\n// Briefly diverge the formatting so this run is not silently re-merged\n// with its still-identical neighbor before the second edit lands.\nrunY.font.bold = true;\nawait context.sync();\n\n// --- Edit 2: append \"123123213\" at the end of the paragraph text -----\nconst contentRange2 = target.getRange(\"Content\");\nconst endRange = contentRange2.getRange(\"End\");\nconst runAppend = endRange.insertText(\"123123213\", \"Replace\");\nrunAppend.font.bold = true;\nawait context.sync();\n\n// --- Restore original (matching) formatting on both new runs ---------\nrunY.font.bold = false;\nrunAppend.font.bold = false;\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# Target change (pr1): the paragraph's single run \"yayaysysys\" becomes\n# three runs with identical rPr (<w:lang w:val=\"en-US\"/>):\n#   \"Y\"  +  \"ayaysysys\"  +  \"123123213\"\n# i.e. the leading \"y\" is capitalized to \"Y\" and \"123123213\" is appended\n# at the end, but the text stays split across three separate <w:r>\n# elements (as real Word leaves them after two distinct edit operations)\n# instead of being coalesced back into one run.\n\n$d = $word.ActiveDocument\n\n# --- Edit 1: capitalize the leading \"y\" -> \"Y\" ------------------------\n$r1 = $d.Content\n$r1.Find.ClearFormatting()\n$r1.Find.Execute(\"y\") | Out-Null\n$r1.Text = \"Y\"\n# Briefly diverge the formatting so this run is not silently re-merged\n# with its still-identical neighbor before the second edit lands.\n$r1.Bold = 1\n\n# --- Edit 2: append \"123123213\" at the end of the paragraph text ------\n$r2 = $d.Range($d.Content.End - 1, $d.Content.End - 1)\n$r2.InsertAfter(\"123123213\")\n$r2.LanguageID = \"en-US\"\n$r2.Bold = 1\n\n# --- Restore original (matching) formatting on both new runs ----------\n$r1.Bold = 0\n$r2.Bold = 0\n"}
